$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.429.66'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '3.535.31'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.95'
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.07'
$ws.Range('E6').Value = '  +2.33%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.594'
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.135'
$ws.Range('E9').Value = '  +8.16%  '
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').Value = '4.140.08'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.77'
$ws.Range('E14').Value = '  +2.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000183'
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('D16').Value = '67.320.85'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').Value = '3.526.25'
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.39'
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.26'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '397.64'
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.03'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.64'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.542'
$ws.Range('E23').Value = '  +2.67%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.31'
$ws.Range('E26').Value = '  +1.62%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.36'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.09'
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '24.16'
$ws.Range('E32').Value = '  +2.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.46'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.63'
$ws.Range('E34').Value = '  +5.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '163.67'
$ws.Range('E35').Value = '  +1.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.902'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.93'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.88'
$ws.Range('E39').Value = '  +2.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0752'
$ws.Range('E40').Value = '  +1.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.83'
$ws.Range('E41').Value = '  +1.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.20'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.64'
$ws.Range('E43').Value = '  +2.52%  '
$ws.Range('D44').Value = '2.813.72'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.01'
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0314'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '343.91'
$ws.Range('E47').Value = '  -3.44%  '
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.78'
$ws.Range('E49').Value = '  +2.30%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.864'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.56'
$ws.Range('E51').Value = '  +1.53%  '
